# Add a new "Tagging" entry to the git-notes workbook (Sheet1), right
# after the existing "My own git server" row (row 38), and move the
# selection/viewport down to reflect the newly added row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# New row 39: Language="Tagging" (A), Topic="Basic" (B), Notes text (C)
$ws.Range("A39").Value = "Tagging"
$ws.Range("B39").Value = "Basic"

$noteText = '# Basic command
`git tag ` //show all tags (both lightweight and annotated)
`git tag v1.1` //lightweight tag HEAD to v1.1
`git tag -a v1.4 -m "my version 1.4"` //build annotated tag v1.4
`git show v1.1` //show the diff between v1.1 and previous commit
# Concept
There are 2 kinds of tags, lightweight and Annotated. Lightweight tag is like a final branch. Annotated tag is a copy in object tree and might go with key and signature and blablabla.'

$ws.Range("C39").Value = $noteText

# Match the wrap-text style used by the other long-form note cells (C37/C38)
$ws.Range("C39").WrapText = $true

# Match row height used for the sibling "My own git server" row above it
$ws.Rows.Item(39).RowHeight = 120

# Reflect the newly inserted row in the saved view: scroll down a couple
# of rows and move the active selection to C40 (one below the new row).
$win = $excel.ActiveWindow
$win.ScrollRow = 34
$win.ScrollColumn = 1

[void]$ws.Range("C40").Select()
